# BALP 1.1.1 and history file updates
#
# - Metadata sheet: bump Version to 1.1.1 and update Date to the new
#   publication timestamp.
# - Elements sheet: clear the stray Constraint(s) text on the root
#   "Extension" row (row 2 / column AI), which the FHIR IG publisher no
#   longer emits for that row.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B3").Value = "1.1.1"
$metadata.Range("B8").Value = "2022-10-21T09:04:31-05:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
